$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9786997437477112
$ws.Range("B1").Value = 1.373613834381104
$ws.Range("C1").Value = 3.679993629455566
$ws.Range("D1").Value = 2.523476362228394
$ws.Range("E1").Value = 0.6842707395553589
